$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data (row 13) continuing the daily log
$ws.Range("A13").Value = 45953
$ws.Range("B13").Value = 5597
$ws.Range("C13").Value = 4312
$ws.Range("D13").Value = 3958
$ws.Range("E13").Value = 272
$ws.Range("F13").Value = 44
$ws.Range("G13").Value = 32
$ws.Range("H13").Value = 6
$ws.Range("I13").Value = 0

# Update selection to H17
$ws.Range("H17").Select()
